# Generate Report for Handoff
# Inserts a new "handoff" row (751d3533-e448-4877-9d0d-898286d1e54e.md) above the
# existing ae7dc5c7-... row on every sheet (Overview / zh-cn / de-de), pushing the
# previous entry down to row 3, and grows the tables/dimensions accordingly.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77b2962b190e50cd10101a81f76f53df3db18564/e2e"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "751d3533-e448-4877-9d0d-898286d1e54e.md"
$ws.Range("B2").Value = "e2e\751d3533-e448-4877-9d0d-898286d1e54e.md"
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-09-02 22:44:35"

$ws.Range("B2").Style = "HyperLink"
$ws.Range("G2").NumberFormat = $dateFmt

foreach ($h in $ws.Hyperlinks) { $h.Delete() }
$ws.Hyperlinks.Add($ws.Range("B2"), "$baseUrl/751d3533-e448-4877-9d0d-898286d1e54e.md", "", "", "e2e\751d3533-e448-4877-9d0d-898286d1e54e.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "$baseUrl/ae7dc5c7-7d5a-430e-8f00-5c467bbbe966.md", "", "", "e2e\ae7dc5c7-7d5a-430e-8f00-5c467bbbe966.md")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "751d3533-e448-4877-9d0d-898286d1e54e.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "751d3533-e448-4877-9d0d-898286d1e54e.d70aecaa4446a88eb369481593db3577492fe356.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-02 22:44:31"
$ws.Range("I2").Value = "'"
$ws.Range("J2").Value = "'"
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = "'"
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = "'"

$ws.Range("A2").Style = "HyperLink"
$ws.Range("H2").NumberFormat = $dateFmt
$ws.Range("K2").NumberFormat = $dateFmt

foreach ($h in $ws.Hyperlinks) { $h.Delete() }
$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/751d3533-e448-4877-9d0d-898286d1e54e.md", "", "", "751d3533-e448-4877-9d0d-898286d1e54e.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/ae7dc5c7-7d5a-430e-8f00-5c467bbbe966.md", "", "", "ae7dc5c7-7d5a-430e-8f00-5c467bbbe966.md")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "751d3533-e448-4877-9d0d-898286d1e54e.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "751d3533-e448-4877-9d0d-898286d1e54e.d70aecaa4446a88eb369481593db3577492fe356.de-de.xlf"
$ws.Range("H2").Value = "2016-09-02 22:44:20"
$ws.Range("I2").Value = "'"
$ws.Range("J2").Value = "'"
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = "'"
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = "'"

$ws.Range("A2").Style = "HyperLink"
$ws.Range("H2").NumberFormat = $dateFmt
$ws.Range("K2").NumberFormat = $dateFmt

foreach ($h in $ws.Hyperlinks) { $h.Delete() }
$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/751d3533-e448-4877-9d0d-898286d1e54e.md", "", "", "751d3533-e448-4877-9d0d-898286d1e54e.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/ae7dc5c7-7d5a-430e-8f00-5c467bbbe966.md", "", "", "ae7dc5c7-7d5a-430e-8f00-5c467bbbe966.md")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

$wb.Save()
